# Update NATMI ligand/receptor TPM-derived values for Ecm1-Itgb4 pair.
# New raw "average expression value" / "total expression value" numbers
# (re-computed from the updated TPM matrix) per cluster, keyed by the
# shared-string id used for that cluster in column A / D
# (20 = ECs, 21 = FAPs, 22 = MuSCs, 23 = Resolving-Mac).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand (Ecm1) average / total expression value per sending cluster
$ligandAvg = @{ 20 = 11.31080466666667; 21 = 101.506841;      22 = 7.161644;          23 = 15.17110633333333 }
$ligandTot = @{ 20 = 33.932414;         21 = 304.520523;      22 = 21.484932;         23 = 45.513319 }

# Updated receptor (Itgb4) average / total expression value per target cluster
$receptorAvg = @{ 20 = 9.771369666666667; 21 = 2.758130333333333; 22 = 1.889356;  23 = 0.7207983333333333 }
$receptorTot = @{ 20 = 29.314109;         21 = 8.274391;          22 = 5.668068;  23 = 2.162395 }

# Receptor-expressing cell count / detection rate per target cluster
$receptorCells = @{ 20 = 3; 21 = 3; 22 = 3; 23 = 3 }
$receptorRate  = @{ 20 = 1; 21 = 1; 22 = 1; 23 = 1 }

# Derived specificity = value / sum(value) across all clusters
$ligandAvgTotal = ($ligandAvg.Values | Measure-Object -Sum).Sum
$ligandTotTotal = ($ligandTot.Values | Measure-Object -Sum).Sum
$receptorAvgTotal = ($receptorAvg.Values | Measure-Object -Sum).Sum
$receptorTotTotal = ($receptorTot.Values | Measure-Object -Sum).Sum

$ligandAvgSpec = @{}
$ligandTotSpec = @{}
foreach ($k in $ligandAvg.Keys) {
    $ligandAvgSpec[$k] = $ligandAvg[$k] / $ligandAvgTotal
    $ligandTotSpec[$k] = $ligandTot[$k] / $ligandTotTotal
}

$receptorAvgSpec = @{}
$receptorTotSpec = @{}
foreach ($k in $receptorAvg.Keys) {
    $receptorAvgSpec[$k] = $receptorAvg[$k] / $receptorAvgTotal
    $receptorTotSpec[$k] = $receptorTot[$k] / $receptorTotTotal
}

# Sheet layout: rows 2-5 sending=ECs(20), 6-9 sending=FAPs(21),
# 10-13 sending=MuSCs(22), 14-17 sending=Resolving-Mac(23); within each
# block of 4 rows the target cluster cycles ECs(20), FAPs(21), MuSCs(22),
# Resolving-Mac(23) (columns A/D, row order matches the worksheet).
$clusters = 20, 21, 22, 23
$row = 2
foreach ($sendCluster in $clusters) {
    foreach ($targetCluster in $clusters) {
        $ws.Range("G$row").Value = $ligandAvg[$sendCluster]
        $ws.Range("H$row").Value = $ligandTot[$sendCluster]
        $ws.Range("I$row").Value = $ligandAvgSpec[$sendCluster]
        $ws.Range("J$row").Value = $ligandTotSpec[$sendCluster]

        $ws.Range("K$row").Value = $receptorCells[$targetCluster]
        $ws.Range("L$row").Value = $receptorRate[$targetCluster]
        $ws.Range("M$row").Value = $receptorAvg[$targetCluster]
        $ws.Range("N$row").Value = $receptorTot[$targetCluster]
        $ws.Range("O$row").Value = $receptorAvgSpec[$targetCluster]
        $ws.Range("P$row").Value = $receptorTotSpec[$targetCluster]

        # Edge weights = ligand metric * receptor metric
        $ws.Range("Q$row").Value = $ligandAvg[$sendCluster] * $receptorAvg[$targetCluster]
        $ws.Range("R$row").Value = $ligandTot[$sendCluster] * $receptorTot[$targetCluster]
        $ws.Range("S$row").Value = $ligandAvgSpec[$sendCluster] * $receptorAvgSpec[$targetCluster]
        $ws.Range("T$row").Value = $ligandTotSpec[$sendCluster] * $receptorTotSpec[$targetCluster]

        $row++
    }
}
